# Updated RAD MRF test case for new Motor Fuel Tax Error message
#
# Refresh the "Date" column (column B) of each payment test-data sheet with
# the timestamps recorded from the latest Katalon RAD test run so the
# BeforePayments.xlsx fixture reflects the current execution rather than the
# stale run captured back in April/October.
#
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Estimated")
$ws.Cells.Item(2, 2).Value = "Wed Oct 16 20:07:25 EDT 2024"
$ws.Cells.Item(3, 2).Value = "Wed Oct 16 20:08:22 EDT 2024"
$ws.Cells.Item(4, 2).Value = "Wed Oct 16 20:09:15 EDT 2024"
$ws.Cells.Item(5, 2).Value = "Wed Oct 16 20:10:08 EDT 2024"
$ws.Cells.Item(6, 2).Value = "Wed Oct 16 20:11:01 EDT 2024"
$ws.Cells.Item(7, 2).Value = "Wed Oct 16 20:11:53 EDT 2024"

$ws = $wb.Worksheets.Item("Existing")
$ws.Cells.Item(2, 2).Value = "Wed Oct 16 20:14:19 EDT 2024"
$ws.Cells.Item(3, 2).Value = "Wed Oct 16 20:15:15 EDT 2024"
$ws.Cells.Item(4, 2).Value = "Wed Oct 16 20:16:08 EDT 2024"
$ws.Cells.Item(5, 2).Value = "Wed Oct 16 20:17:02 EDT 2024"
$ws.Cells.Item(6, 2).Value = "Wed Oct 16 20:17:56 EDT 2024"
$ws.Cells.Item(7, 2).Value = "Wed Oct 16 20:18:48 EDT 2024"
$ws.Cells.Item(8, 2).Value = "Wed Oct 16 20:19:42 EDT 2024"
$ws.Cells.Item(9, 2).Value = "Wed Oct 16 20:20:36 EDT 2024"
$ws.Cells.Item(10, 2).Value = "Wed Oct 16 20:21:30 EDT 2024"
$ws.Cells.Item(11, 2).Value = "Wed Oct 16 20:22:23 EDT 2024"
$ws.Cells.Item(12, 2).Value = "Wed Oct 16 20:23:17 EDT 2024"
$ws.Cells.Item(13, 2).Value = "Wed Oct 16 20:24:10 EDT 2024"
$ws.Cells.Item(14, 2).Value = "Wed Oct 16 20:25:04 EDT 2024"
$ws.Cells.Item(15, 2).Value = "Wed Oct 16 20:25:57 EDT 2024"
$ws.Cells.Item(16, 2).Value = "Wed Oct 16 20:26:50 EDT 2024"
$ws.Cells.Item(17, 2).Value = "Wed Oct 16 20:27:44 EDT 2024"
$ws.Cells.Item(18, 2).Value = "Wed Oct 16 20:28:37 EDT 2024"
$ws.Cells.Item(19, 2).Value = "Wed Oct 16 20:29:30 EDT 2024"

$ws = $wb.Worksheets.Item("Extension")
$ws.Cells.Item(2, 2).Value = "Wed Oct 16 21:27:47 EDT 2024"
$ws.Cells.Item(3, 2).Value = "Wed Oct 16 21:28:43 EDT 2024"
$ws.Cells.Item(4, 2).Value = "Wed Oct 16 21:29:36 EDT 2024"
$ws.Cells.Item(5, 2).Value = "Wed Oct 16 21:30:28 EDT 2024"
$ws.Cells.Item(6, 2).Value = "Wed Oct 16 21:31:19 EDT 2024"
$ws.Cells.Item(7, 2).Value = "Wed Oct 16 21:32:11 EDT 2024"

$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Cells.Item(2, 2).Value = "Wed Oct 16 21:33:03 EDT 2024"
$ws.Cells.Item(3, 2).Value = "Wed Oct 16 21:33:54 EDT 2024"
$ws.Cells.Item(4, 2).Value = "Wed Oct 16 21:34:46 EDT 2024"
$ws.Cells.Item(5, 2).Value = "Wed Oct 16 21:35:37 EDT 2024"
$ws.Cells.Item(6, 2).Value = "Wed Oct 16 21:36:29 EDT 2024"
$ws.Cells.Item(7, 2).Value = "Wed Oct 16 21:37:20 EDT 2024"
$ws.Cells.Item(8, 2).Value = "Wed Oct 16 21:38:12 EDT 2024"
$ws.Cells.Item(9, 2).Value = "Wed Oct 16 21:39:04 EDT 2024"
$ws.Cells.Item(10, 2).Value = "Wed Oct 16 21:39:55 EDT 2024"
$ws.Cells.Item(11, 2).Value = "Wed Oct 16 21:40:47 EDT 2024"
$ws.Cells.Item(12, 2).Value = "Wed Oct 16 21:41:38 EDT 2024"
$ws.Cells.Item(13, 2).Value = "Wed Oct 16 21:42:29 EDT 2024"
$ws.Cells.Item(14, 2).Value = "Wed Oct 16 21:43:21 EDT 2024"
$ws.Cells.Item(15, 2).Value = "Wed Oct 16 21:44:12 EDT 2024"
$ws.Cells.Item(16, 2).Value = "Wed Oct 16 21:45:03 EDT 2024"
$ws.Cells.Item(17, 2).Value = "Wed Oct 16 21:45:55 EDT 2024"
$ws.Cells.Item(18, 2).Value = "Wed Oct 16 21:46:48 EDT 2024"
$ws.Cells.Item(19, 2).Value = "Wed Oct 16 21:47:41 EDT 2024"
$ws.Cells.Item(20, 2).Value = "Wed Oct 16 21:48:34 EDT 2024"
$ws.Cells.Item(21, 2).Value = "Wed Oct 16 21:49:26 EDT 2024"
$ws.Cells.Item(22, 2).Value = "Wed Oct 16 21:50:18 EDT 2024"
$ws.Cells.Item(23, 2).Value = "Wed Oct 16 21:51:11 EDT 2024"
$ws.Cells.Item(24, 2).Value = "Wed Oct 16 21:52:03 EDT 2024"
$ws.Cells.Item(25, 2).Value = "Wed Oct 16 21:52:56 EDT 2024"
$ws.Cells.Item(26, 2).Value = "Wed Oct 16 21:53:49 EDT 2024"
$ws.Cells.Item(27, 2).Value = "Wed Oct 16 21:54:42 EDT 2024"
$ws.Cells.Item(28, 2).Value = "Wed Oct 16 21:55:35 EDT 2024"
$ws.Cells.Item(29, 2).Value = "Wed Oct 16 21:56:29 EDT 2024"
$ws.Cells.Item(30, 2).Value = "Wed Oct 16 21:57:22 EDT 2024"
$ws.Cells.Item(31, 2).Value = "Wed Oct 16 21:58:15 EDT 2024"
$ws.Cells.Item(32, 2).Value = "Wed Oct 16 21:59:08 EDT 2024"
$ws.Cells.Item(33, 2).Value = "Wed Oct 16 22:00:00 EDT 2024"
$ws.Cells.Item(34, 2).Value = "Wed Oct 16 22:00:53 EDT 2024"
$ws.Cells.Item(35, 2).Value = "Wed Oct 16 22:01:47 EDT 2024"
$ws.Cells.Item(36, 2).Value = "Wed Oct 16 22:02:39 EDT 2024"
$ws.Cells.Item(37, 2).Value = "Wed Oct 16 22:03:32 EDT 2024"
$ws.Cells.Item(38, 2).Value = "Wed Oct 16 22:04:26 EDT 2024"
$ws.Cells.Item(39, 2).Value = "Wed Oct 16 22:05:20 EDT 2024"
$ws.Cells.Item(40, 2).Value = "Wed Oct 16 22:06:12 EDT 2024"
$ws.Cells.Item(41, 2).Value = "Wed Oct 16 22:07:05 EDT 2024"
$ws.Cells.Item(42, 2).Value = "Wed Oct 16 22:07:59 EDT 2024"
$ws.Cells.Item(43, 2).Value = "Wed Oct 16 22:08:51 EDT 2024"
$ws.Cells.Item(44, 2).Value = "Wed Oct 16 22:09:44 EDT 2024"
$ws.Cells.Item(45, 2).Value = "Wed Oct 16 22:10:36 EDT 2024"
$ws.Cells.Item(46, 2).Value = "Wed Oct 16 22:11:29 EDT 2024"
$ws.Cells.Item(47, 2).Value = "Wed Oct 16 22:12:22 EDT 2024"
$ws.Cells.Item(48, 2).Value = "Wed Oct 16 22:13:15 EDT 2024"
$ws.Cells.Item(49, 2).Value = "Wed Oct 16 22:14:08 EDT 2024"
$ws.Cells.Item(50, 2).Value = "Wed Oct 16 22:15:01 EDT 2024"
$ws.Cells.Item(51, 2).Value = "Wed Oct 16 22:15:54 EDT 2024"
$ws.Cells.Item(52, 2).Value = "Wed Oct 16 22:16:47 EDT 2024"

$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Cells.Item(2, 2).Value = "Wed Oct 16 22:19:22 EDT 2024"
$ws.Cells.Item(3, 2).Value = "Wed Oct 16 22:20:11 EDT 2024"
$ws.Cells.Item(4, 2).Value = "Wed Oct 16 22:21:00 EDT 2024"
$ws.Cells.Item(5, 2).Value = "Wed Oct 16 22:21:49 EDT 2024"
$ws.Cells.Item(6, 2).Value = "Wed Oct 16 22:22:39 EDT 2024"
$ws.Cells.Item(7, 2).Value = "Wed Oct 16 22:23:28 EDT 2024"
$ws.Cells.Item(8, 2).Value = "Wed Oct 16 22:24:17 EDT 2024"
$ws.Cells.Item(9, 2).Value = "Wed Oct 16 22:25:07 EDT 2024"

$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Cells.Item(2, 2).Value = "Wed Oct 16 22:25:56 EDT 2024"
$ws.Cells.Item(3, 2).Value = "Wed Oct 16 22:26:51 EDT 2024"
$ws.Cells.Item(4, 2).Value = "Wed Oct 16 22:27:47 EDT 2024"
$ws.Cells.Item(5, 2).Value = "Wed Oct 16 22:28:42 EDT 2024"
$ws.Cells.Item(6, 2).Value = "Wed Oct 16 22:29:38 EDT 2024"

$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Cells.Item(2, 2).Value = "Wed Oct 16 22:17:40 EDT 2024"
$ws.Cells.Item(3, 2).Value = "Wed Oct 16 22:18:30 EDT 2024"
